$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing t_value_* columns (AN:AY) - these metrics no longer
# apply to this general-version sheet.
$ws.Range("AN1:AY4").EntireColumn.Delete()

# Remove the old "Run_3" row (row 4) entirely.
$ws.Rows(4).Delete()

# Update the remaining header labels in columns AL/AM.
$ws.Cells.Item(1, 38).Value = "t_value_YX_C"
$ws.Cells.Item(1, 39).Value = "t_value_pH_UL"

# Row 3 (formerly "Run_2") becomes "Run_3" with refreshed metric values.
$ws.Cells.Item(3, 1).Value = "Run_3"
$ws.Cells.Item(3, 2).Value = 0.3019000248478971
$ws.Cells.Item(3, 3).Value = 6.270094823160237

$ws.Cells.Item(3, 21).Value = 0.1767723542181233
$ws.Cells.Item(3, 22).Value = 0.1801977329232616
$ws.Cells.Item(3, 23).Value = 21.53459878428799
$ws.Cells.Item(3, 24).Value = 1.288233347656097
$ws.Cells.Item(3, 25).Value = 0.5653948624759768
$ws.Cells.Item(3, 26).Value = 27.81498117706033
$ws.Cells.Item(3, 27).Value = 0.05908581434424913
$ws.Cells.Item(3, 28).Value = 0.4657913127924631
$ws.Cells.Item(3, 29).Value = 5.922211725760993
$ws.Cells.Item(3, 30).Value = 0.0780941576563732
$ws.Cells.Item(3, 31).Value = 0.09885336412199139
$ws.Cells.Item(3, 32).Value = 0.9373952689530096
$ws.Cells.Item(3, 33).Value = -28.79586067230095
$ws.Cells.Item(3, 34).Value = -27.21234173384484
$ws.Cells.Item(3, 35).Value = 0.6519937603508629
$ws.Cells.Item(3, 36).Value = 0.0936729665876723
$ws.Cells.Item(3, 37).Value = 14.05229673901558
$ws.Cells.Item(3, 38).Value = 650.6571407306653
$ws.Cells.Item(3, 39).Value = 2350319.829006001
